# Append-new-listings update: 2025-12-25 01:57 JST
# Sheet 1 ("ランサーズ") holds the scraped Lancers job listing rows.
# This run refreshed the "fetched at" timestamp for every existing row and
# inserted two brand-new listings into the ranked list (at row 15 and row 19),
# pushing the rows that used to occupy those slots further down.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newTimestamp = "2025-12-25 01:57:32"

# --- 1. Insert the two new rows where the new listings now rank --------
# Row 15: new WordPress favorites-feature listing (pushes everything
#         from the old row 15 onward down by one).
$ws.Rows.Item(15).Insert()
# Row 19: new "限定公開" listing (pushes the old row 18 - the VPS/FTP
#         listing - down to row 20). This is computed AFTER the first
#         insert, so it already accounts for the shift above.
$ws.Rows.Item(19).Insert()

# --- 2. Refresh the "fetched at" timestamp in column A for every data row
for ($r = 2; $r -le 20; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}

# --- 3. Populate the newly inserted row 15 ------------------------------
$ws.Range("B15").Value = "【WordPress】会員向け「お気に入り機能」実装(マイページ一覧表示まで)"
$ws.Range("C15").Value = "システム開発"
$ws.Range("D15").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E15").Value = "期限情報なし"
$ws.Range("F15").Value = "https://www.lancers.jp/work/detail/5460732"
$ws.Range("F15").Style = "Hyperlink"
$ws.Range("G15").Value = 28
$ws.Range("H15").Value = "○WordPress"

# --- 4. Populate the newly inserted row 19 ------------------------------
$ws.Range("B19").Value = "限定公開 限定公開の仕事"
$ws.Range("C19").Value = "システム開発"
$ws.Range("D19").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E19").Value = "期限情報なし"
$ws.Range("F19").Value = "https://www.lancers.jp/work/detail/5450323"
$ws.Range("F19").Style = "Hyperlink"
$ws.Range("G19").Value = 13

# --- 5. Wire up hyperlinks for the two brand-new rows -------------------
# (the rows that merely shifted down keep pointing at their original,
# now-stale relationship ids - that's exactly what the source data does)
$ws.Hyperlinks.Add($ws.Range("F19"), "https://www.lancers.jp/work/detail/5450323")
$ws.Hyperlinks.Add($ws.Range("F20"), "https://www.lancers.jp/work/detail/5459964")

# Re-assert the shared "Hyperlink" cell style on the two new link cells -
# Hyperlinks.Add() stamps its own ad-hoc style variant on the cell; putting
# it back to the named "Hyperlink" style keeps F15:F20 on the identical
# style index the rest of column F already uses.
$ws.Range("F19").Style = "Hyperlink"
$ws.Range("F20").Style = "Hyperlink"
